$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = 0.82457518239957484
$ws.Range("R1").Value = 0.95165478894049993
$ws.Range("AB1").Value = 0.92835711672778565
$ws.Range("AG1").Value = 0.87857064535627472
$ws.Range("AK1").Value = 0.80815949527109088
$ws.Range("G2").Value = 0.84038372785456383
$ws.Range("O3").Value = 0.91163922004679454
$ws.Range("AX3").Value = 0.67764609578088419
$ws.Range("I4").Value = 0.71085114913570413
$ws.Range("AK4").Value = 0.78905594985915462
$ws.Range("BM4").Value = 0.78765458520582809
$ws.Range("C5").Value = 0.87890192375165155
$ws.Range("F5").Value = 0.74781636869155166
$ws.Range("AZ5").Value = 0.96180690328870089
$ws.Range("BH6").Value = 0.94940332826969542
$ws.Range("H7").Value = 0.98288141123843376
$ws.Range("BB7").Value = 0.76253761278865728
$ws.Range("AE8").Value = 0.63293196181144706
$ws.Range("P9").Value = 0.79591540846234587
$ws.Range("I10").Value = 0.98094935249477055
$ws.Range("U10").Value = 0.79033662383687697
$ws.Range("AZ10").Value = 0.99776071290238866
$ws.Range("I11").Value = 0.54952613435773401
$ws.Range("BA11").Value = 0.85915874449924168
$ws.Range("BG11").Value = 0.58797026174475475
$ws.Range("M12").Value = 0.93644195258038732
$ws.Range("AH13").Value = 0.94231387175579573
$ws.Range("A14").Value = 0.84508111148025544
$ws.Range("L14").Value = 0.99543164108131599
$ws.Range("H15").Value = 0.69469823420101617
$ws.Range("Z15").Value = 0.70010691240903444
$ws.Range("X16").Value = 0.86095216675967412
$ws.Range("AL16").Value = 0.89322858812677075
$ws.Range("P17").Value = 0.66472836737092811
$ws.Range("S17").Value = 0.80991403901495462
$ws.Range("T17").Value = 0.71784950176734696
$ws.Range("AU18").Value = 0.82031973526688295
$ws.Range("AV19").Value = 0.70545819884801575
$ws.Range("AK21").Value = 0.59403406612390064
$ws.Range("BC21").Value = 0.80896804346147433
$ws.Range("BI22").Value = 0.97132398272717335
$ws.Range("U23").Value = 0.69664062421945228
$ws.Range("BH23").Value = 0.95536144213237484
$ws.Range("AE24").Value = 0.93641636035234677
$ws.Range("AJ24").Value = 0.97984892088967213
$ws.Range("AF25").Value = 0.97089257745458335
$ws.Range("BP25").Value = 0.82105522185639024
$ws.Range("BG26").Value = 0.93137444715586848
$ws.Range("BE27").Value = 0.83756203341115487
$ws.Range("S28").Value = 0.65317033337004582
$ws.Range("AA29").Value = 0.68974511531890392
$ws.Range("AL29").Value = 0.68551797710971951
$ws.Range("BD29").Value = 0.98990017832572752
$ws.Range("BE29").Value = 0.8774687743440478
$ws.Range("AF30").Value = 0.69614633241654533
$ws.Range("AS30").Value = 0.86629538302205344
$ws.Range("AC31").Value = 0.94948644699603313
$ws.Range("B32").Value = 0.63460423083447515
$ws.Range("K32").Value = 0.62804651825040458
$ws.Range("U32").Value = 0.85247417288815075
$ws.Range("Z32").Value = 0.69897974449537337
$ws.Range("AI33").Value = 0.90758270293542631
$ws.Range("AV33").Value = 0.91589429344529938
$ws.Range("BO34").Value = 0.88807031419399052
$ws.Range("N35").Value = 0.98023032139934574
$ws.Range("V35").Value = 0.72056312412182266
$ws.Range("AH36").Value = 0.99953580854119994
$ws.Range("AI36").Value = 0.74808838599705618
$ws.Range("BP38").Value = 0.92577633501814638
$ws.Range("B39").Value = 0.92542940240141391
$ws.Range("V39").Value = 0.7214577955350614
$ws.Range("AI39").Value = 0.8052876965074085
$ws.Range("AB40").Value = 0.99645318373614367
$ws.Range("AH40").Value = 0.6945334988988312
$ws.Range("AZ40").Value = 0.96364179143561934
$ws.Range("D41").Value = 0.97129321323421669
$ws.Range("Y41").Value = 0.8328844194955427
$ws.Range("AT41").Value = 0.90559163668805376
$ws.Range("AM42").Value = 0.79908669224985795
$ws.Range("AV42").Value = 0.80834794679972477
$ws.Range("BD42").Value = 0.8734567898252269
$ws.Range("F43").Value = 0.90839485315536073
$ws.Range("R43").Value = 0.55233582083180544
$ws.Range("W43").Value = 0.76904390753192875
$ws.Range("AR43").Value = 0.80496736221562348
$ws.Range("AJ44").Value = 0.91219500714499069
$ws.Range("AO44").Value = 0.94950587406496845
$ws.Range("BH44").Value = 0.70437070108865507
$ws.Range("P45").Value = 0.78798596958759504
$ws.Range("AB45").Value = 0.84572294835804041
$ws.Range("T46").Value = 0.85844692304560866
$ws.Range("A47").Value = 0.92268783182265934
$ws.Range("M48").Value = 0.86446211315092492
$ws.Range("AS48").Value = 0.9122338742146413
$ws.Range("G49").Value = 0.86467417129555357
$ws.Range("AI49").Value = 0.90833664810493941
$ws.Range("BN49").Value = 0.99029109191273013
$ws.Range("P50").Value = 0.76755012624316721
$ws.Range("BH50").Value = 0.98808001568417447
$ws.Range("X51").Value = 0.68952866932952084
$ws.Range("AC51").Value = 0.86411412503237417
$ws.Range("BF51").Value = 0.98859331949123996
$ws.Range("BG51").Value = 0.89701113987582382
$ws.Range("B52").Value = 0.6502169357733939
$ws.Range("T52").Value = 0.95406578230323058
$ws.Range("AA52").Value = 0.85367431233771485
$ws.Range("AX53").Value = 0.65932608583407104
$ws.Range("BK53").Value = 0.96074487213587667
$ws.Range("BL53").Value = 0.87866143527047269
$ws.Range("L54").Value = 0.92701258394531894
$ws.Range("Z54").Value = 0.70377855153826707
$ws.Range("H55").Value = 0.7833586403172732
$ws.Range("A56").Value = 0.66058247138005455
$ws.Range("P56").Value = 0.76992351267690928
$ws.Range("Y58").Value = 0.88557596239509229
$ws.Range("AN59").Value = 0.92641690022578305
$ws.Range("BH59").Value = 0.99687176487610973
$ws.Range("BK61").Value = 0.93015360902199951
$ws.Range("BN61").Value = 0.78232612711037031
$ws.Range("T62").Value = 0.72000008925330827
$ws.Range("BK62").Value = 0.93340805963783291
$ws.Range("BP62").Value = 0.68605431504423997
$ws.Range("F63").Value = 0.82561280603923382
$ws.Range("AK63").Value = 0.75809762260598612
$ws.Range("AN63").Value = 0.95830223021895644
$ws.Range("AR63").Value = 0.77365135459423506
$ws.Range("BD63").Value = 0.93581938839847756
$ws.Range("S64").Value = 0.68288043912005181
$ws.Range("AN64").Value = 0.81366523784002476
$ws.Range("BE64").Value = 0.8526953983065555
$ws.Range("A65").Value = 0.6217602196782277
$ws.Range("W66").Value = 0.66784553038080019
$ws.Range("I67").Value = 0.88090956476533355
$ws.Range("J67").Value = 0.95107946012329458
$ws.Range("V67").Value = 0.86832884537665522
$ws.Range("W67").Value = 0.94965037013812514
